# Update the GitHub repository URL text on slide 3 (shape "CustomShape 1", id 230):
#   "/MSD_R_course_TT2024_2 " (trailing space included) -> "/MSD_R_course_MT2024_2"
# This also merges the trailing single-space run into the edited run and removes it,
# matching the authored edit (TT2024 -> MT2024, trailing stray space deleted).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

$full = $tr.Text
$old = "/MSD_R_course_TT2024_2 "
$new = "/MSD_R_course_MT2024_2"

$idx = $full.IndexOf($old)
if ($idx -ge 0) {
    $startPos = $idx + 1
    $len = $old.Length
    $c = $tr.Characters($startPos, $len)
    $c.Text = $new
}
